# "so sanh dang nv" — add the "Qualities Capacity Level" numbers (column E)
# to the Input_Asset sheet, size the "Usage Logs" column to fit its content,
# best-fit the Input_Employee qualities columns, and restore the sheet
# selections/active sheet that were captured when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# --- Input_Asset: fill in column E (Qualities Capacity Level) ---------------
$wsAsset = $wb.Worksheets.Item("Input_Asset")

$wsAsset.Range("E2").Value = 1
$wsAsset.Range("E3").Value = 2
$wsAsset.Range("E4").Value = 1
$wsAsset.Range("E5").Value = 4
$wsAsset.Range("E6").Value = 2
$wsAsset.Range("E7").Value = 1

# Column F ("Usage Logs") holds a long JSON blob — widen it to fit.
$wsAsset.Columns("F").AutoFit() | Out-Null

# --- Input_Task: no data changes, just restore the saved selection ---------
$wsTask = $wb.Worksheets.Item("Input_Task")
$wsTask.Range("A1:H15").Select() | Out-Null

# --- Input_Employee: best-fit the qualities columns, restore selection -----
$wsEmployee = $wb.Worksheets.Item("Input_Employee")
$wsEmployee.Columns("C:M").AutoFit() | Out-Null
$wsEmployee.Range("H15").Select() | Out-Null

# --- Leave Input_Asset as the active sheet with E8 selected -----------------
$wsAsset.Activate() | Out-Null
$wsAsset.Range("E8").Select() | Out-Null
